$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.905.24"
$ws.Range("E2").Value = "  +1.61%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.620.55"
$ws.Range("E3").Value = "  +5.06%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.30"
$ws.Range("E5").Value = "  +4.76%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "320.09"
$ws.Range("E6").Value = "  +0.48%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +0.30%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.37"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11 - Chainlink
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.76"
$ws.Range("E11").Value = "  -1.24%  "

# Row 12 - Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0805"
$ws.Range("E12").Value = "  +0.60%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.61%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.18"
$ws.Range("E14").Value = "  +1.66%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.033.92"
$ws.Range("E15").Value = "  +5.23%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.611.55"
$ws.Range("E16").Value = "  +5.60%  "

# Row 17 - Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("E17").Value = "  +1.68%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.813.55"
$ws.Range("E18").Value = "  +1.67%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +0.26%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  +1.28%  "

# Row 21 - ImmutableX
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.87"
$ws.Range("E21").Value = "  -1.97%  "

# Row 22 - ShibaInu
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0938"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.44"
$ws.Range("E23").Value = "  -3.98%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.92"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.91%  "

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.18"
$ws.Range("E26").Value = "  +1.90%  "

# Row 27 - Dai
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.05%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  +5.09%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +0.62%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.35"
$ws.Range("E30").Value = "  +1.99%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  -1.40%  "

# Row 32 - OKB
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.26"
$ws.Range("E32").Value = "  +0.48%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +2.71%  "

# Row 34 - FirstDigitalUSD
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35 - Celestia
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.02"
$ws.Range("E35").Value = "  -1.49%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0788"
$ws.Range("E36").Value = "  +2.30%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.90"
$ws.Range("E37").Value = "  +8.78%  "

# Row 38 - ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.01"
$ws.Range("E38").Value = "  +3.43%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +7.45%  "

# Row 40 - Monero
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.18"
$ws.Range("E40").Value = "  +4.30%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -0.51%  "

# Row 42 - EnergySwap
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.17"
$ws.Range("E42").Value = "  +2.64%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  -0.41%  "

# Row 44 - VeChain
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0312"
$ws.Range("E44").Value = "  +4.42%  "

# Row 45 - Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.077.95"
$ws.Range("E45").Value = "  +4.45%  "

# Row 46 - NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("E46").Value = "  +2.26%  "

# Row 47 - Stacks
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +9.04%  "

# Row 48 - ApeXProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.20"
$ws.Range("E48").Value = "  +4.60%  "

# Row 49 - was FraxShare, now RocketPoolETH (new coin takes this slot)
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.884.80"
$ws.Range("E49").Value = "  +5.44%  "

# Row 50 - was THORChain, now FraxShare (shifted down, THORChain removed from list)
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.84"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51 - MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.66"
$ws.Range("E51").Value = "  +4.68%  "
